# Generate Report for Handoff
# Adds a new tracked file (fb2474c2-bbf8-4669-9c40-e0ea0ca1bba9.md) as row 7
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$guidName   = "fb2474c2-bbf8-4669-9c40-e0ea0ca1bba9.md"
$guidPath   = "e2e\fb2474c2-bbf8-4669-9c40-e0ea0ca1bba9.md"
$zhXlf      = "fb2474c2-bbf8-4669-9c40-e0ea0ca1bba9.7d3f986024dee741048e9954fc2b33b51d4b6806.zh-cn.xlf"
$deXlf      = "fb2474c2-bbf8-4669-9c40-e0ea0ca1bba9.7d3f986024dee741048e9954fc2b33b51d4b6806.de-de.xlf"

$hoDateOverview = "2016-09-06 10:07:43"
$hoDateZh       = "2016-09-06 10:07:33"
$hoDateDe       = "2016-09-06 10:07:43"
$epoch          = "0001-01-01 00:00:00"

$srcUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0000000000000000000000000000000000000000/e2e/fb2474c2-bbf8-4669-9c40-e0ea0ca1bba9.md"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0000000000000000000000000000000000000000/e2e/fb2474c2-bbf8-4669-9c40-e0ea0ca1bba9.md"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0000000000000000000000000000000000000000/e2e/fb2474c2-bbf8-4669-9c40-e0ea0ca1bba9.md"

function Set-BlankCell($range) {
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview  (columns A-G)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A7").Value2 = $guidName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), $srcUrl, "", "", $guidPath) | Out-Null
$wsOverview.Range("C7").Value2 = ".md"
Set-BlankCell $wsOverview.Range("D7")
$wsOverview.Range("E7").Value2 = "Ready for handoff"
$wsOverview.Range("F7").Value2 = "Ready for handoff"
$wsOverview.Range("G7").Value2 = $hoDateOverview
$wsOverview.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn  (columns A-P)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A7"), $zhUrl, "", "", $guidName) | Out-Null
$wsZh.Range("B7").Value2 = ".md"
$wsZh.Range("C7").Value2 = "Ready for handoff"
$wsZh.Range("D7").Value2 = "e2e"
$wsZh.Range("E7").Value2 = "ht"
$wsZh.Range("F7").Value2 = "False"
$wsZh.Range("G7").Value2 = $zhXlf
$wsZh.Range("H7").Value2 = $hoDateZh
$wsZh.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-BlankCell $wsZh.Range("I7")
Set-BlankCell $wsZh.Range("J7")
$wsZh.Range("K7").Value2 = $epoch
$wsZh.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-BlankCell $wsZh.Range("L7")
$wsZh.Range("M7").Value2 = "True"
Set-BlankCell $wsZh.Range("N7")
$wsZh.Range("O7").Value2 = "False"
Set-BlankCell $wsZh.Range("P7")

# ---------------------------------------------------------------------------
# Sheet 3: de-de  (columns A-P)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A7"), $deUrl, "", "", $guidName) | Out-Null
$wsDe.Range("B7").Value2 = ".md"
$wsDe.Range("C7").Value2 = "Ready for handoff"
$wsDe.Range("D7").Value2 = "e2e"
$wsDe.Range("E7").Value2 = "ht"
$wsDe.Range("F7").Value2 = "False"
$wsDe.Range("G7").Value2 = $deXlf
$wsDe.Range("H7").Value2 = $hoDateDe
$wsDe.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-BlankCell $wsDe.Range("I7")
Set-BlankCell $wsDe.Range("J7")
$wsDe.Range("K7").Value2 = $epoch
$wsDe.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
Set-BlankCell $wsDe.Range("L7")
$wsDe.Range("M7").Value2 = "True"
Set-BlankCell $wsDe.Range("N7")
$wsDe.Range("O7").Value2 = "False"
Set-BlankCell $wsDe.Range("P7")

$wb.Save()
